$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Games": append two newly-played games (rows 40 and 41) that were
# previously listed as upcoming games in the "Next" sheet.
# ---------------------------------------------------------------------------
$games = $wb.Worksheets.Item("Games")

# Row 40 -> Game 39, played 2024-01-13 (serial 45306) vs LAL
$games.Cells.Item(40, 1).Value = 39
$games.Cells.Item(40, 2).Value = 45306
$games.Cells.Item(40, 2).NumberFormat = "YYYY-MM-DD"
$games.Cells.Item(40, 3).Value = -1
$games.Cells.Item(40, 4).Value = 105
$games.Cells.Item(40, 5).Value = 97
$games.Cells.Item(40, 6).Value = 0.495
$games.Cells.Item(40, 7).Value = 9.800000000000001
$games.Cells.Item(40, 8).Value = 26.5
$games.Cells.Item(40, 9).Value = 0.104
$games.Cells.Item(40, 10).Value = 108.3
$games.Cells.Item(40, 11).Value = "LAL"
$games.Cells.Item(40, 12).Value = 112
$games.Cells.Item(40, 13).Value = 0.571
$games.Cells.Item(40, 14).Value = 15.6
$games.Cells.Item(40, 15).Value = 22.9
$games.Cells.Item(40, 16).Value = 0.19
$games.Cells.Item(40, 17).Value = 115.5
$games.Cells.Item(40, 18).Value = 0
$games.Cells.Item(40, 19).Value = 0

# Row 41 -> Game 40, played 2024-01-14 (serial 45307) vs LAC
$games.Cells.Item(41, 1).Value = 40
$games.Cells.Item(41, 2).Value = 45307
$games.Cells.Item(41, 2).NumberFormat = "YYYY-MM-DD"
$games.Cells.Item(41, 3).Value = -2
$games.Cells.Item(41, 4).Value = 117
$games.Cells.Item(41, 5).Value = 91.90000000000001
$games.Cells.Item(41, 6).Value = 0.602
$games.Cells.Item(41, 7).Value = 11.6
$games.Cells.Item(41, 8).Value = 25
$games.Cells.Item(41, 9).Value = 0.205
$games.Cells.Item(41, 10).Value = 127.4
$games.Cells.Item(41, 11).Value = "LAC"
$games.Cells.Item(41, 12).Value = 128
$games.Cells.Item(41, 13).Value = 0.655
$games.Cells.Item(41, 14).Value = 8.5
$games.Cells.Item(41, 15).Value = 31.6
$games.Cells.Item(41, 16).Value = 0.161
$games.Cells.Item(41, 17).Value = 139.3
$games.Cells.Item(41, 18).Value = 0
$games.Cells.Item(41, 19).Value = 0

# ---------------------------------------------------------------------------
# Sheet "Next": those two games (LAL @ 45306 and LAC @ 45307) already
# happened, so remove them from the upcoming-games list; everything below
# shifts up by two rows.
# ---------------------------------------------------------------------------
$next = $wb.Worksheets.Item("Next")
$next.Range("A2:C3").EntireRow.Delete()
